$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change in condition labelling: update the "neutral"/"emotional" labels
# to "control"/"positive" labels for the toy condition study (row 7).
$ws.Range("E7").Value = "no-toy / control: 0.2"
$ws.Range("D7").Value = "toy / positive: 0.4"

# Update the active selection to match the edited cell area.
$ws.Range("D7").Select()
